$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 155. Excel shifts the existing rows
# 155..269 down to 156..270 (carrying formatting, including the date
# number-format on column D), and the sheet's used range / dimension grows
# to A1:R270 automatically.
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new weekly price record.
$ws.Cells.Item(155, 1).Value = 3
$ws.Cells.Item(155, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(155, 3).Value = "Coquimbo"
$ws.Cells.Item(155, 4).Value = 44574
$ws.Cells.Item(155, 5).Value = 5
$ws.Cells.Item(155, 6).Value = 100112009
$ws.Cells.Item(155, 7).Value = "Acelga"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 235
$ws.Cells.Item(155, 11).Value = 2200
$ws.Cells.Item(155, 12).Value = 2500
$ws.Cells.Item(155, 13).Value = 2347
$ws.Cells.Item(155, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(155, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(155, 16).Value = 391
$ws.Cells.Item(155, 17).Value = 6
$ws.Cells.Item(155, 18).Value = "Hortaliza"
